# Generate Report for Handoff
# Updates the localization-status report after a new handoff run:
#  - Overview sheet: "Latest HO Xliff Generate Date" timestamps refreshed
#    for the files that were just handed off again (rows 7,8,10,11,13,14)
#  - de-de sheet: "Latest Handoff Datetime" shares the same refreshed
#    timestamp as the Overview sheet for those rows
#  - zh-cn sheet: "Latest Handoff Datetime" refreshed with its own
#    (slightly earlier) timestamp for those rows
#  - zh-cn / de-de sheets: "Priority" set to "ht" for the handed-off rows

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 10, 11, 13, 14)

$overviewTimestamp = "2016-09-06 18:28:31"
$zhCnTimestamp = "2016-09-06 18:28:22"

$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = $overviewTimestamp
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 8).Value = $zhCnTimestamp
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 8).Value = $overviewTimestamp
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
}
